$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Semestre ideal" value in column B/C (EQD-6 -> EQD-7)
$ws.Range("B9").Value = "EQD-7,EQN-11"
$ws.Range("C9").Value = "EQD-7,EQN-11"

# Add the new "Requisitos" row (row 26) with the new weak-requirement entry,
# copying the formatting used by the existing requirement row above it (row 25)
$ws.Range("B25:C25").Copy($ws.Range("B26:C26"))

$newReq = "LOQ4082 -  Corrosão  (Requisito fraco)`n"
$ws.Range("B26").Value = $newReq
$ws.Range("C26").Value = $newReq

$ws.Rows.Item(26).RowHeight = 30
